# This edit swaps the data contained in row 4 and row 5 of the sheet
# (everything except the header row remains structurally the same column
# set per-row, but the two records change places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New row 4 (previously the data that lived in row 5)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = 60493331
$ws.Range("B4").Value = 106964
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 220299
$ws.Range("F4").Value = "Svinrot"
$ws.Range("G4").Value = "Scorzonera humilis"
$ws.Range("H4").Value = "L."
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("P4").Value = "Svalboviken, Sågtorp, Srm"
$ws.Range("Q4").Value = 575959.7454739227
$ws.Range("R4").Value = 6561005.400135099
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Södermanland"
$ws.Range("U4").Value = "Eskilstuna"
$ws.Range("V4").Value = "Södermanland"
$ws.Range("W4").Value = "Näshulta"
$ws.Range("Y4").Value = "'2016-07-03"
$ws.Range("Z4").Value = "00:00"
$ws.Range("AA4").Value = "'2016-07-03"
$ws.Range("AB4").Value = "00:00"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AI4").Value = "Skogsbryn mot väg"
$ws.Range("AW4").Value = "Håkan Gustafson"
$ws.Range("AX4").Value = "Håkan Gustafson"

# ---------------------------------------------------------------------
# New row 5 (previously the data that lived in row 4)
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 87015259
$ws.Range("B5").Value = 44335
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 102021
$ws.Range("F5").Value = "Mindre bastardsvärmare"
$ws.Range("G5").Value = "Zygaena viciae"
$ws.Range("H5").Value = "(Denis & Schiffermüller, 1775)"
$ws.Range("I5").Value = "'1"
$ws.Range("J5").Value = "ex."
$ws.Range("K5").Value = "imago/adult"
$ws.Range("L5").ClearContents()
$ws.Range("M5").Value = "födosökande"
$ws.Range("P5").Value = "Nytorp 2, Eskilstuna, Srm"
$ws.Range("Q5").Value = 575983.1057489648
$ws.Range("R5").Value = 6561147.345532569
$ws.Range("S5").Value = 25
$ws.Range("T5").Value = "Södermanland"
$ws.Range("U5").Value = "Eskilstuna"
$ws.Range("V5").Value = "Södermanland"
$ws.Range("W5").Value = "Näshulta"
$ws.Range("Y5").Value = "'2020-07-21"
$ws.Range("Z5").Value = "17:00"
$ws.Range("AA5").Value = "'2020-07-21"
$ws.Range("AB5").Value = "17:00"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AI5").ClearContents()
$ws.Range("AW5").Value = "Thomas Holmgren"
$ws.Range("AX5").Value = "Thomas Holmgren"
